# "Implement a new runner"
#
# The original deck had several text boxes whose label was split across two
# <a:r> runs (e.g. "Parameters and " + "constraints"). This script merges
# each of those runs back into a single run, and repositions a connector
# line + callout group on slide 2 that were moved as part of the same edit.

$p = $ppt.ActivePresentation

function Merge-ShapeText {
    param($shape, [string]$text)

    # Reading TextRange.Text already returns the concatenation of every run
    # in the paragraph, so if we set it directly to the same concatenation
    # the runtime sees "no change" and never rewrites the underlying runs.
    # Setting a throw-away value first forces the real write to happen.
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "~"
    $tr.Text = $text
}

# ---------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

Merge-ShapeText $s1.Shapes.Item(4).GroupItems.Item(1) "Parameters and constraints"
Merge-ShapeText $s1.Shapes.Item(8) "Parameters and constraints"
Merge-ShapeText $s1.Shapes.Item(16).GroupItems.Item(1) "Parameters and constraints"
Merge-ShapeText $s1.Shapes.Item(24) "Decoded Covering Array"
Merge-ShapeText $s1.Shapes.Item(33).GroupItems.Item(1) "Parameters and constraints"
Merge-ShapeText $s1.Shapes.Item(37) "Decoded Covering Array"
Merge-ShapeText $s1.Shapes.Item(41) "Parameters and constraints"
Merge-ShapeText $s1.Shapes.Item(55) "Decoded Covering Array"
Merge-ShapeText $s1.Shapes.Item(57) "Decoded Covering Array"

# ---------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

Merge-ShapeText $s2.Shapes.Item(19) "Covering Array"
Merge-ShapeText $s2.Shapes.Item(20) "Covering Array"
Merge-ShapeText $s2.Shapes.Item(21) "Covering Array"
Merge-ShapeText $s2.Shapes.Item(22) "Covering Array"
Merge-ShapeText $s2.Shapes.Item(45) "Covering Array"
Merge-ShapeText $s2.Shapes.Item(46).GroupItems.Item(1) "Parameters and constraints"
Merge-ShapeText $s2.Shapes.Item(55) "Decoded Covering Array"

# Reposition / resize the dotted "Straight Connector 130" line.
# PowerPoint's object model reports Left/Top/Width/Height in points
# (1 pt = 12700 EMU); the target values below are the diff's EMU values
# converted to points.
$conn = $s2.Shapes.Item(59)
$conn.Left = 6050941 / 12700
$conn.Top = 3599826 / 12700
$conn.Width = 722744 / 12700
$conn.Height = 1363069 / 12700

# Move the "Group 144" callout group (size stays the same).
$grp = $s2.Shapes.Item(60)
$grp.Left = 4991486 / 12700
$grp.Top = 2959249 / 12700
